$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (existing row 11 and below shift down to 12, etc.)
$ws.Rows.Item(11).Insert()

# Bring over the cell formatting (number formats / alignment / wrap) from row 9,
# which carries the same style pattern we need for the new row 11.
$ws.Range("A9:B9").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("C9").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("D9:E9").Copy()
$ws.Range("D11:E11").PasteSpecial(-4122)

$ws.Range("F9:I9").Copy()
$ws.Range("F11:I11").PasteSpecial(-4122)

# Fill in the data for the new row (Entrevista - 2 Corte - blind person interview)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Entrevista"
$ws.Range("C11").Value = 43531
$ws.Range("D11").Value = 0.41666666666666669
$ws.Range("E11").Value = 0.47916666666666669
$ws.Range("F11").Value = 10
$ws.Range("G11").Value = 80
$ws.Range("H11").Value = "2 Corte"
$ws.Range("I11").Value = "Se llamó a la persona ciega y se preguntó por dudas y otras cuestiones relacionadas con el proceso para personas con necesidades especiales"

# Match the taller row height this new entry needs.
$ws.Rows.Item(11).RowHeight = 72

# Update the view state: scrolled down a bit, with the new last cell selected.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("I12").Select()
